$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = 278184
$ws.Range("E10").Value = 1751974896
$ws.Range("C17").Value = 134740
$ws.Range("E17").Value = 296788602
$ws.Range("C19").Value = 108913
$ws.Range("D19").Value = 13558
$ws.Range("E19").Value = 344646966
$ws.Range("C65").Value = 61053
$ws.Range("E65").Value = 334085242
$ws.Range("C81").Value = 26157
$ws.Range("E81").Value = 165046014
$ws.Range("C85").Value = 10747
$ws.Range("E85").Value = 47050687
$ws.Range("C97").Value = 98505
$ws.Range("E97").Value = 307060090
$ws.Range("C104").Value = 22088
$ws.Range("E104").Value = 84759332
$ws.Range("C110").Value = 16865
$ws.Range("E110").Value = 25928635
$ws.Range("C115").Value = 17536
$ws.Range("E115").Value = 38564635
$ws.Range("C117").Value = 19691
$ws.Range("E117").Value = 56394544
$ws.Range("C122").Value = 9687
$ws.Range("E122").Value = 31911349
$ws.Range("C132").Value = 6664
$ws.Range("E132").Value = 13449261
$ws.Range("C134").Value = 5664
$ws.Range("E134").Value = 17026661
$ws.Range("C150").Value = 95008
$ws.Range("E150").Value = 278791733
$ws.Range("C152").Value = 126038
$ws.Range("E152").Value = 715740995
$ws.Range("C164").Value = 50560
$ws.Range("E164").Value = 168354513
$ws.Range("C168").Value = 284895
$ws.Range("E168").Value = 1207328625
$ws.Range("C169").Value = 562551
$ws.Range("E169").Value = 1284022282
$ws.Range("C170").Value = 367207
$ws.Range("E170").Value = 2842859211
$ws.Range("C173").Value = 54378
$ws.Range("E173").Value = 151815965
$ws.Range("C174").Value = 357124
$ws.Range("E174").Value = 1015769599
$ws.Range("C175").Value = 125481
$ws.Range("E175").Value = 809697991
$ws.Range("C179").Value = 235623
$ws.Range("E179").Value = 811591522
$ws.Range("C180").Value = 141438
$ws.Range("E180").Value = 339801579
$ws.Range("C186").Value = 21927
$ws.Range("E186").Value = 39904691
$ws.Range("C188").Value = 19685
$ws.Range("E188").Value = 65970654
$ws.Range("C189").Value = 2008
$ws.Range("E189").Value = 8001252
$ws.Range("C192").Value = 7455
$ws.Range("E192").Value = 16969326
$ws.Range("C193").Value = 5342
$ws.Range("E193").Value = 27691588
$ws.Range("C196").Value = 7393
$ws.Range("E196").Value = 20466364
$ws.Range("C197").Value = 6965
$ws.Range("E197").Value = 9390903
$ws.Range("C199").Value = 4154
$ws.Range("E199").Value = 9011715
$ws.Range("C203").Value = 13090
$ws.Range("E203").Value = 32919097
$ws.Range("C208").Value = 1533
$ws.Range("E208").Value = 3221269
$ws.Range("C213").Value = 3628
$ws.Range("E213").Value = 11063625
$ws.Range("C220").Value = 4711
$ws.Range("E220").Value = 11655990
$ws.Range("C257").Value = 182545
$ws.Range("E257").Value = 1063734456
$ws.Range("C275").Value = 39298
$ws.Range("E275").Value = 100474879
$ws.Range("C276").Value = 216629
$ws.Range("E276").Value = 1209877432
$ws.Range("C284").Value = 48484
$ws.Range("E284").Value = 65577519
$ws.Range("C287").Value = 8568
$ws.Range("E287").Value = 24938099
$ws.Range("C307").Value = 39617
$ws.Range("E307").Value = 95225547
$ws.Range("C311").Value = 190852
$ws.Range("E311").Value = 586331627
$ws.Range("C317").Value = 103575
$ws.Range("E317").Value = 302856654

$wb.Save()
